$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 42.75280866666666
$ws.Range("H2").Value = 128.258426
$ws.Range("I2").Value = 0.8529286054750734
$ws.Range("J2").Value = 0.8529286054750735
$ws.Range("M2").Value = 25.11140833333333
$ws.Range("N2").Value = 75.334225
$ws.Range("O2").Value = 0.7431105026796001
$ws.Range("P2").Value = 0.7431105026796001
$ws.Range("Q2").Value = 1073.583235825539
$ws.Range("R2").Value = 9662.24912242985
$ws.Range("S2").Value = 0.6338202047643922
$ws.Range("T2").Value = 0.6338202047643923
$ws.Range("G3").Value = 42.75280866666666
$ws.Range("H3").Value = 128.258426
$ws.Range("I3").Value = 0.8529286054750734
$ws.Range("J3").Value = 0.8529286054750735
$ws.Range("O3").Value = 0.1596166092346045
$ws.Range("P3").Value = 0.1596166092346045
$ws.Range("Q3").Value = 230.6005839180982
$ws.Range("R3").Value = 2075.405255262884
$ws.Range("S3").Value = 0.1361415719251309
$ws.Range("T3").Value = 0.136141571925131
$ws.Range("G4").Value = 42.75280866666666
$ws.Range("H4").Value = 128.258426
$ws.Range("I4").Value = 0.8529286054750734
$ws.Range("J4").Value = 0.8529286054750735
$ws.Range("N4").Value = 9.861222
$ws.Range("O4").Value = 0.09727288808579543
$ws.Range("P4").Value = 0.09727288808579541
$ws.Range("Q4").Value = 140.5316457951747
$ws.Range("R4").Value = 1264.784812156572
$ws.Range("S4").Value = 0.08296682878555038
$ws.Range("T4").Value = 0.08296682878555038
$ws.Range("I5").Value = 0.04642608686423023
$ws.Range("J5").Value = 0.04642608686423023
$ws.Range("M5").Value = 25.11140833333333
$ws.Range("N5").Value = 75.334225
$ws.Range("O5").Value = 0.7431105026796001
$ws.Range("P5").Value = 0.7431105026796001
$ws.Range("Q5").Value = 58.4366244049889
$ws.Range("R5").Value = 525.9296196449001
$ws.Range("S5").Value = 0.03449971274712491
$ws.Range("T5").Value = 0.03449971274712491
$ws.Range("I6").Value = 0.04642608686423023
$ws.Range("J6").Value = 0.04642608686423023
$ws.Range("O6").Value = 0.1596166092346045
$ws.Range("P6").Value = 0.1596166092346045
$ws.Range("S6").Value = 0.00741037456529964
$ws.Range("T6").Value = 0.007410374565299641
$ws.Range("I7").Value = 0.04642608686423023
$ws.Range("J7").Value = 0.04642608686423023
$ws.Range("N7").Value = 9.861222
$ws.Range("O7").Value = 0.09727288808579543
$ws.Range("P7").Value = 0.09727288808579541
$ws.Range("Q7").Value = 7.649332374338668
$ws.Range("R7").Value = 68.843991369048
$ws.Range("S7").Value = 0.004515999551805684
$ws.Range("T7").Value = 0.004515999551805683
$ws.Range("G8").Value = 5.044817999999999
$ws.Range("I8").Value = 0.1006453076606963
$ws.Range("J8").Value = 0.1006453076606963
$ws.Range("M8").Value = 25.11140833333333
$ws.Range("N8").Value = 75.334225
$ws.Range("O8").Value = 0.7431105026796001
$ws.Range("P8").Value = 0.7431105026796001
$ws.Range("Q8").Value = 126.68248476535
$ws.Range("R8").Value = 1140.14236288815
$ws.Range("S8").Value = 0.07479058516808303
$ws.Range("T8").Value = 0.07479058516808303
$ws.Range("G9").Value = 5.044817999999999
$ws.Range("I9").Value = 0.1006453076606963
$ws.Range("J9").Value = 0.1006453076606963
$ws.Range("O9").Value = 0.1596166092346045
$ws.Range("P9").Value = 0.1596166092346045
$ws.Range("Q9").Value = 27.21079650300399
$ws.Range("S9").Value = 0.0160646627441739
$ws.Range("T9").Value = 0.01606466274417391
$ws.Range("G10").Value = 5.044817999999999
$ws.Range("I10").Value = 0.1006453076606963
$ws.Range("J10").Value = 0.1006453076606963
$ws.Range("N10").Value = 9.861222
$ws.Range("O10").Value = 0.09727288808579543
$ws.Range("P10").Value = 0.09727288808579541
$ws.Range("S10").Value = 0.00979005974843936
$ws.Range("T10").Value = 0.009790059748439358
